# Excel : Ajout type graphique
#
# Adds a 4th column "Graphique-resultats" to the "Tableau2" list object on
# the "Infos" sheet, fills in the sample value "Camenbert" for row 2, and
# adds a list-type data validation on D2 offering "Camenbert, Diagramme,
# Jauge" as the chart-type choices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Infos")
$lo = $ws.ListObjects.Item(1)

# Add the new table column (auto-extends the table ref, autoFilter and the
# worksheet dimension from C2 to D2)
$newCol = $lo.ListColumns.Add()

# Header cell (D1): text + centered (no wrap) alignment, like the other
# header cells in that row
$hdrCell = $ws.Range("D1")
$hdrCell.Value2 = "Graphique-resultats"
$hdrCell.HorizontalAlignment = -4108
$hdrCell.VerticalAlignment = -4108
$hdrCell.WrapText = $false

# Data cell (D2): sample value + same centered alignment
$dataCell = $ws.Range("D2")
$dataCell.Value2 = "Camenbert"
$dataCell.HorizontalAlignment = -4108
$dataCell.VerticalAlignment = -4108
$dataCell.WrapText = $false

# Data validation list on D2 for the chart type
$dataCell.Validation.Add(3, $null, $null, '"Camenbert, Diagramme, Jauge"')

# Restore the selection Excel left behind after inserting the column
$ws.Range("D7").Select()

Write-Output ("Tableau2 now spans " + $lo.Range.Address())
